# PROS-14717 LionJP Adjacency KPI
# Rename the kpi_name value in B2 from PRODUCT_GROUP_ADJACENCY_IN_WHOLE_STORE
# to ADJACENCY_PRODUCT_GROUP_IN_SCENE_TYPE, and move the active selection to B9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "ADJACENCY_PRODUCT_GROUP_IN_SCENE_TYPE"

$ws.Range("B9").Select()
